$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indexed")

# Delete the "iL" column (column G), shifting columns H:K left to G:J.
$ws.Range("G1").EntireColumn.Delete()
